# EPBDS-13605 Performance Improvement: The expression in the return cell is
# executed, despite the fact that the corresponding rule was not matched.
#
# Update the "times" expectation row (rows 36, 81, 122, 167 on the "Testing"
# sheet) so that each indicator (times/times2/times3/result) is expected to
# be incremented only once instead of multiple times, reflecting that
# short-circuited (non-matched) conditions must no longer invoke their
# expressions.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Testing")

$ws2.Range("D36:H36").Value   = 1
$ws2.Range("D81:H81").Value   = 1
$ws2.Range("D122:H122").Value = 1
$ws2.Range("D167:H167").Value = 1

# Re-establish the merged header ranges so they line up with the refreshed
# layout (this also matches the reordering Excel performs on the
# <mergeCells> list after re-merging a touched range).
$m1 = $ws2.Range("B129:E129")
$m1.UnMerge()
$m1.Merge()

$m2 = $ws2.Range("B152:C152")
$m2.UnMerge()
$m2.Merge()

# Move the active selection/view to the cell that was last edited, and drop
# the previous scroll anchor.
[void]$ws2.Activate()
[void]$ws2.Range("D36").Select()
